$wb = $excel.ActiveWorkbook

# "Handed back" status text shared by both locale sheets.
$handedBackStatus = "Handed back: in sync with en-US"

# New "Latest Handback DateTime" values (row 2 = the 0c168e65... file,
# row 3 = the c8686c97... file) per locale sheet.
$handbackTimes = @{
    "zh-cn" = "2016-02-24 06:49:56"
    "de-de" = "2016-02-24 06:50:24"
}

foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Build a lookup of existing hyperlink addresses keyed by the cell's
    # address ($A$2, $C$2, ...) so the new "Latest Target File" (E) /
    # "Latest Handback File" (F) columns can reuse the very same targets
    # as the existing "Source File Name" (A) / "Latest Handoff File" (C)
    # hyperlinks.
    $linkByCell = @{}
    foreach ($hl in $ws.Hyperlinks) {
        $linkByCell[$hl.Range.Address()] = $hl.Address
    }

    $aAddr2 = $linkByCell['$A$2']
    $cAddr2 = $linkByCell['$C$2']
    $aAddr3 = $linkByCell['$A$3']
    $cAddr3 = $linkByCell['$C$3']

    $aDisp2 = $ws.Range("A2").Value2
    $cDisp2 = $ws.Range("C2").Value2
    $aDisp3 = $ws.Range("A3").Value2
    $cDisp3 = $ws.Range("C3").Value2

    # Row 2 (0c168e65-...md handoff pair): mark handed back, populate the
    # Latest Target File / Latest Handback File columns, and stamp the real
    # handback datetime.
    $ws.Range("B2").Value2 = $handedBackStatus
    $ws.Range("G2").Value2 = $handbackTimes[$sheetName]

    $ws.Range("E2").Value2 = $aDisp2
    $ws.Hyperlinks.Add($ws.Range("E2"), $aAddr2, "", "", $aDisp2) | Out-Null
    $ws.Range("E2").Font.Underline = $true
    $ws.Range("E2").Font.Color = 15570276

    $ws.Range("F2").Value2 = $cDisp2
    $ws.Hyperlinks.Add($ws.Range("F2"), $cAddr2, "", "", $cDisp2) | Out-Null
    $ws.Range("F2").Font.Underline = $true
    $ws.Range("F2").Font.Color = 15570276

    # Row 3 (c8686c97-...md handoff pair): same treatment.
    $ws.Range("B3").Value2 = $handedBackStatus
    $ws.Range("G3").Value2 = $handbackTimes[$sheetName]

    $ws.Range("E3").Value2 = $aDisp3
    $ws.Hyperlinks.Add($ws.Range("E3"), $aAddr3, "", "", $aDisp3) | Out-Null
    $ws.Range("E3").Font.Underline = $true
    $ws.Range("E3").Font.Color = 15570276

    $ws.Range("F3").Value2 = $cDisp3
    $ws.Hyperlinks.Add($ws.Range("F3"), $cAddr3, "", "", $cDisp3) | Out-Null
    $ws.Range("F3").Font.Underline = $true
    $ws.Range("F3").Font.Color = 15570276
}
